# Adds three new LAF test cases (TC_07..TC_09) to the "LAF" sheet and
# flips the "Run" flag (column C) of the pre-existing Google-search related
# test cases (rows 5-7) from "Yes" to "No" now that the new cases are the
# active ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Flip existing rows 5-7 "Run" column from Yes -> No -------------------
$ws.Range("C5").Value2 = "No"
$ws.Range("C6").Value2 = "No"
$ws.Range("C7").Value2 = "No"

# --- Populate the three new rows -------------------------------------------
# Values are entered column-by-column (A, B, C, F, then BX) so that new
# shared-string entries are created in the same order the workbook expects.

# Column A - TestScenario
$ws.Range("A8").Value2  = "TC_07_Validate_LAF_Title"
$ws.Range("A9").Value2  = "TC_08_Validate_JoinNow_Button"
$ws.Range("A10").Value2 = "TC_09_Validate_JoinNow_ClickButton"

# Column B - Module
$ws.Range("B8").Value2  = "Membership"
$ws.Range("B9").Value2  = "Membership"
$ws.Range("B10").Value2 = "Membership"

# Column C - Run
$ws.Range("C8").Value2  = "Yes"
$ws.Range("C9").Value2  = "Yes"
$ws.Range("C10").Value2 = "Yes"

# Column F - TextMessage
$ws.Range("F8").Value2  = "CheckLAF_Title"
$ws.Range("F9").Value2  = "Check_LAF_JoinNowHeaderButton"
$ws.Range("F10").Value2 = "Click_JoinNowHeaderButton"

# Column BX - Text_input (only needed on the title-validation row)
$ws.Range("BX8").Value2 = "LA Fitness | Gym and Fitness Club | Join Today"

# --- Formatting fixups -------------------------------------------------
# C7 (now a real "No" data row) and the three new C8:C10 cells need the
# same left/top aligned style already used by C5/C6.
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C7:C10").PasteSpecial(-4122) | Out-Null

# BX5 no longer needs its distinct fill style; match the plain/default
# formatting used elsewhere in the row (e.g. A5).
$ws.Range("A5").Copy() | Out-Null
$ws.Range("BX5").PasteSpecial(-4122) | Out-Null

# --- Selection / view bookkeeping ------------------------------------------
$ws.Range("F17").Select() | Out-Null

Write-Output "LAF test cases added"
